$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2752.054
$ws.Range("I98").Value = 2937.9355
$ws.Range("J98").Value = 1791.6666
$ws.Range("K98").Value = 2937.9355
$ws.Range("L98").Value = 1791.6666
$ws.Range("M98").Value = -1439.9355
$ws.Range("N98").Value = -4787.6666
$ws.Range("H122").Value = 2752.054
$ws.Range("I122").Value = 2937.9355
$ws.Range("J122").Value = 1791.6666
$ws.Range("K122").Value = 8813.806500000001
$ws.Range("L122").Value = 5374.9998
$ws.Range("M122").Value = -6363.806500000001
$ws.Range("N122").Value = -10274.9998
$ws.Range("H132").Value = 2826.0952
$ws.Range("I132").Value = 2975.1082
$ws.Range("J132").Value = 1723.4
$ws.Range("K132").Value = 8925.3246
$ws.Range("L132").Value = 5170.200000000001
$ws.Range("M132").Value = -6395.3246
$ws.Range("N132").Value = -10230.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 841835.9
$ws.Range("I2").Value = 1320.7391
$ws.Range("J2").Value = 2452823.2
$ws.Range("K2").Value = 1320.7391
$ws.Range("L2").Value = 2452823.2
$ws.Range("M2").Value = -1207.7391
$ws.Range("N2").Value = -2453049.2
$ws.Range("H32").Value = 8560.105
$ws.Range("I32").Value = 4592.027
$ws.Range("K32").Value = 4592.027
$ws.Range("M32").Value = -4305.027
$ws.Range("H74").Value = 13640417
$ws.Range("I74").Value = 18750704
$ws.Range("K74").Value = 18750704
$ws.Range("M74").Value = -18749830
$ws.Range("H77").Value = 13640417
$ws.Range("I77").Value = 18750704
$ws.Range("K77").Value = 93753520
$ws.Range("M77").Value = -93749152
$ws.Range("H116").Value = 841835.9
$ws.Range("I116").Value = 1320.7391
$ws.Range("J116").Value = 2452823.2
$ws.Range("K116").Value = 1320.7391
$ws.Range("L116").Value = 2452823.2
$ws.Range("M116").Value = 973.2609
$ws.Range("N116").Value = -2457411.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 841835.9
$ws.Range("I3").Value = 1320.7391
$ws.Range("J3").Value = 2452823.2
$ws.Range("K3").Value = 1320.7391
$ws.Range("L3").Value = 2452823.2
$ws.Range("M3").Value = -1206.7391
$ws.Range("N3").Value = -2453051.2
$ws.Range("H107").Value = 428710.5
$ws.Range("I107").Value = 565131.2
$ws.Range("J107").Value = 2395.75
$ws.Range("K107").Value = 565131.2
$ws.Range("L107").Value = 2395.75
$ws.Range("M107").Value = -563211.2
$ws.Range("N107").Value = -6235.75

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2276.24
$ws.Range("I31").Value = 1814.6666
$ws.Range("J31").Value = 2968.6
$ws.Range("K31").Value = 1814.6666
$ws.Range("L31").Value = 2968.6
$ws.Range("M31").Value = -1519.6666
$ws.Range("N31").Value = -3558.6
$ws.Range("H34").Value = 2276.24
$ws.Range("I34").Value = 1814.6666
$ws.Range("J34").Value = 2968.6
$ws.Range("K34").Value = 1814.6666
$ws.Range("L34").Value = 2968.6
$ws.Range("M34").Value = -1612.6666
$ws.Range("N34").Value = -3372.6
$ws.Range("H86").Value = 3726.9333
$ws.Range("I86").Value = 3734.3333
$ws.Range("J86").Value = 3722
$ws.Range("K86").Value = 3734.3333
$ws.Range("L86").Value = 3722
$ws.Range("M86").Value = -2611.3333
$ws.Range("N86").Value = -5968
$ws.Range("H89").Value = 3726.9333
$ws.Range("I89").Value = 3734.3333
$ws.Range("J89").Value = 3722
$ws.Range("K89").Value = 18671.6665
$ws.Range("L89").Value = 18610
$ws.Range("M89").Value = -13055.6665
$ws.Range("N89").Value = -29842
$ws.Range("H94").Value = 1158.25
$ws.Range("I94").Value = 653.9
$ws.Range("J94").Value = 1998.8334
$ws.Range("K94").Value = 653.9
$ws.Range("L94").Value = 1998.8334
$ws.Range("M94").Value = -202.9
$ws.Range("N94").Value = -2900.8334
$ws.Range("H103").Value = 10000
$ws.Range("I103").Value = 10000
$ws.Range("K103").Value = 10000
$ws.Range("M103").Value = -8828
$ws.Range("H132").Value = 3623.35
$ws.Range("I132").Value = 3488.0715
$ws.Range("J132").Value = 3939
$ws.Range("K132").Value = 10464.2145
$ws.Range("L132").Value = 11817
$ws.Range("M132").Value = -7934.2145
$ws.Range("N132").Value = -16877
$ws.Range("H134").Value = 2356.4285
$ws.Range("I134").Value = 2154.8147
$ws.Range("K134").Value = 6464.4441
$ws.Range("M134").Value = -3929.4441

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 25000524
$ws.Range("I5").Value = 381.33334
$ws.Range("J5").Value = 45455184
$ws.Range("K5").Value = 1144.00002
$ws.Range("L5").Value = 136365552
$ws.Range("M5").Value = -1032.00002
$ws.Range("N5").Value = -136365776
$ws.Range("H34").Value = 2329.652
$ws.Range("I34").Value = 266.5
$ws.Range("J34").Value = 3057.8235
$ws.Range("K34").Value = 799.5
$ws.Range("L34").Value = 9173.470499999999
$ws.Range("M34").Value = -715.5
$ws.Range("N34").Value = -9341.470499999999
$ws.Range("H80").Value = 10347852
$ws.Range("I80").Value = 18107616
$ws.Range("J80").Value = 1500
$ws.Range("K80").Value = 54322848
$ws.Range("L80").Value = 4500
$ws.Range("M80").Value = -54321912
$ws.Range("N80").Value = -6372
$ws.Range("H83").Value = 10347852
$ws.Range("I83").Value = 18107616
$ws.Range("J83").Value = 1500
$ws.Range("K83").Value = 162968544
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -162963864
$ws.Range("N83").Value = -22860
$ws.Range("H92").Value = 375
$ws.Range("I92").Value = 250
$ws.Range("K92").Value = 750
$ws.Range("M92").Value = 498
$ws.Range("H107").Value = 117934.82
$ws.Range("I107").Value = 250281
$ws.Range("J107").Value = 77212.92
$ws.Range("K107").Value = 750843
$ws.Range("L107").Value = 231638.76
$ws.Range("M107").Value = -748923
$ws.Range("N107").Value = -235478.76
$ws.Range("H113").Value = 501.0811
$ws.Range("I113").Value = 499.17648
$ws.Range("J113").Value = 502.7
$ws.Range("K113").Value = 1497.52944
$ws.Range("L113").Value = 1508.1
$ws.Range("M113").Value = 672.47056
$ws.Range("N113").Value = -5848.1
$ws.Range("H131").Value = 746.38
$ws.Range("J131").Value = 781.5412
$ws.Range("L131").Value = 2344.6236
$ws.Range("N131").Value = -12424.6236
$ws.Range("H132").Value = 47620420
$ws.Range("I132").Value = 76923860
$ws.Range("J132").Value = 2335
$ws.Range("K132").Value = 692314740
$ws.Range("L132").Value = 21015
$ws.Range("M132").Value = -692312210
$ws.Range("N132").Value = -26075
$ws.Range("H135").Value = 25000524
$ws.Range("I135").Value = 381.33334
$ws.Range("J135").Value = 45455184
$ws.Range("K135").Value = 3432.00006
$ws.Range("L135").Value = 409096656
$ws.Range("M135").Value = -897.0000600000003
$ws.Range("N135").Value = -409101726

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2007
$ws.Range("I132").Value = 1846
$ws.Range("J132").Value = 2122
$ws.Range("K132").Value = 5538
$ws.Range("L132").Value = 6366
$ws.Range("M132").Value = -3008
$ws.Range("N132").Value = -11426

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1510.4375
$ws.Range("I61").Value = 1465.5714
$ws.Range("J61").Value = 1824.5
$ws.Range("K61").Value = 1465.5714
$ws.Range("L61").Value = 1824.5
$ws.Range("M61").Value = -1263.5714
$ws.Range("N61").Value = -2228.5
$ws.Range("H113").Value = 1510.4375
$ws.Range("I113").Value = 1465.5714
$ws.Range("J113").Value = 1824.5
$ws.Range("K113").Value = 1465.5714
$ws.Range("L113").Value = 1824.5
$ws.Range("M113").Value = 704.4286
$ws.Range("N113").Value = -6164.5
$ws.Range("H122").Value = 1917.9584
$ws.Range("I122").Value = 1626
$ws.Range("J122").Value = 1994.7894
$ws.Range("K122").Value = 4878
$ws.Range("L122").Value = 5984.3682
$ws.Range("M122").Value = -2428
$ws.Range("N122").Value = -10884.3682

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1179.8049
$ws.Range("I126").Value = 1228.6428
$ws.Range("J126").Value = 1074.6154
$ws.Range("K126").Value = 3685.9284
$ws.Range("L126").Value = 3223.8462
$ws.Range("M126").Value = -1215.9284
$ws.Range("N126").Value = -8163.8462
